$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 236: open/high/low/close change from 177680000000000 to 177772000000000 ---
$ws.Range("C236").Value = 177772000000000
$ws.Range("D236").Value = 177772000000000
$ws.Range("E236").Value = 177772000000000
$ws.Range("F236").Value = 177772000000000

# --- Append new rows 237-239, carrying formatting from row 236 (date style, borders, etc.) ---
$ws.Range("A236:G236").Copy($ws.Range("A237:G237"))
$ws.Range("A236:G236").Copy($ws.Range("A238:G238"))
$ws.Range("A236:G236").Copy($ws.Range("A239:G239"))

# Row 237 : 2023-07-01
$ws.Range("A237").Value = 45108.41666666666
$ws.Range("B237").Value = "ECONOMICS:IQM2"
$ws.Range("C237").Value = 176788000000000
$ws.Range("D237").Value = 176788000000000
$ws.Range("E237").Value = 176788000000000
$ws.Range("F237").Value = 176788000000000
$ws.Range("G237").Value = 0

# Row 238 : 2023-08-01
$ws.Range("A238").Value = 45139.41666666666
$ws.Range("B238").Value = "ECONOMICS:IQM2"
$ws.Range("C238").Value = 174322000000000
$ws.Range("D238").Value = 174322000000000
$ws.Range("E238").Value = 174322000000000
$ws.Range("F238").Value = 174322000000000
$ws.Range("G238").Value = 0

# Row 239 : 2023-09-01
$ws.Range("A239").Value = 45170.41666666666
$ws.Range("B239").Value = "ECONOMICS:IQM2"
$ws.Range("C239").Value = 173950000000000
$ws.Range("D239").Value = 173950000000000
$ws.Range("E239").Value = 173950000000000
$ws.Range("F239").Value = 173950000000000
$ws.Range("G239").Value = 0
